# Regenerated "想去人数" (interest count) figures for the 杭州-漫展信息 workbook.
# Touches the F column (numeric interest counts) across the 展览, 演出 and
# 全部类型 sheets; 本地生活 has no data rows and is untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 852
$ws.Range("F5").Value = 1176
$ws.Range("F6").Value = 13
$ws.Range("F7").Value = 3841
$ws.Range("F8").Value = 2543
$ws.Range("F10").Value = 2401
$ws.Range("F14").Value = 1627
$ws.Range("F17").Value = 95
$ws.Range("F18").Value = 303
$ws.Range("F20").Value = 49
$ws.Range("F21").Value = 264
$ws.Range("F23").Value = 437
$ws.Range("F26").Value = 482
$ws.Range("F27").Value = 668
$ws.Range("F28").Value = 85
$ws.Range("F30").Value = 369
$ws.Range("F31").Value = 38
$ws.Range("F32").Value = 1609
$ws.Range("F33").Value = 867
$ws.Range("F34").Value = 35
$ws.Range("F35").Value = 8
$ws.Range("F36").Value = 928
$ws.Range("F37").Value = 1952
$ws.Range("F39").Value = 515
$ws.Range("F42").Value = 591
$ws.Range("F43").Value = 1237
$ws.Range("F44").Value = 32
$ws.Range("F46").Value = 410

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 11

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 852
$ws.Range("F3").Value = 1176
$ws.Range("F5").Value = 13
$ws.Range("F6").Value = 3841
$ws.Range("F7").Value = 2543
$ws.Range("F8").Value = 2401
$ws.Range("F10").Value = 1627
$ws.Range("F14").Value = 95
$ws.Range("F15").Value = 303
$ws.Range("F17").Value = 49
$ws.Range("F18").Value = 264
$ws.Range("F20").Value = 437
$ws.Range("F23").Value = 482
$ws.Range("F24").Value = 668
$ws.Range("F25").Value = 85
$ws.Range("F30").Value = 369
$ws.Range("F31").Value = 38
$ws.Range("F32").Value = 1609
$ws.Range("F33").Value = 867
$ws.Range("F34").Value = 35
$ws.Range("F36").Value = 928
$ws.Range("F37").Value = 1952
$ws.Range("F42").Value = 515
$ws.Range("F45").Value = 591
$ws.Range("F46").Value = 1237
$ws.Range("F47").Value = 32
$ws.Range("F48").Value = 410
$ws.Range("F49").Value = 11
